$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.415.81'
$ws.Range("E2").Value = '  +0.18%  '
$ws.Range("D3").Value = '1.849.83'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9998'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '240.82'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.17%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6304'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.05%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07709'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.27%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2943'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.43%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.51'
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07749'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.62%  '
$ws.Range("D12").Value = '1.853.18'
$ws.Range("E12").Value = '  +0.19%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.025'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.63%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.00001088'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +8.78%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6802'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.39%  '
$ws.Range("D17").Value = '2.102.31'
$ws.Range("E17").Value = '  -0.11%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.154'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.63%  '
$ws.Range("D19").Value = '29.437.94'
$ws.Range("E19").Value = '  +0.09%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '229.54'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.95%  '
$ws.Range("E21").Value = '  +0.32%  '
$ws.Range("E22").Value = '  +0.06%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.453'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.12%  '
$ws.Range("E24").Value = '  +0.07%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '157.50'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.37%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1389'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.51%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.366'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.19%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.68'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.24%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.468'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.24%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.310'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +4.71%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.05748'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.28%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.113'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.17%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.052'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.76%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.852'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.46%  '
$ws.Range("E35").Value = '  +0.52%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7094'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.58%  '
$ws.Range("E37").Value = '  -0.30%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.777'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.02%  '
$ws.Range("B39").Value = 'Maker'
$ws.Range("C39").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D39").Value = '1.229.61'
$ws.Range("E39").Value = '  -2.43%  '
$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01799'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.81%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.478'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +4.10%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9130'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.15%  '
$ws.Range("E43").Value = '  +0.08%  '
$ws.Range("D44").Value = '2.011.20'
$ws.Range("E44").Value = '  -0.13%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '101.75'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.62%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '66.31'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.30%  '
$ws.Range("E47").Value = '  +3.76%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.154'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.46%  '
$ws.Range("E49").Value = '  -0.58%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.056'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.19%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.687'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.32%  '
